$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the header row: "MaxStart" -> "MaxStar"
$ws.Range("C2").Value = "MaxStar"

# Re-write B2's value (still "IsOpen") so the shared-string table gets
# rebuilt the same way Excel would when the nearby cell text is edited
$ws.Range("B2").Value = "IsOpen"

# Move the active cell/selection to G5
$ws.Range("G5").Select()

$wb.Save()
